# epexspot_prices.xlsx - automatic data refresh
# Sheet "Prix Spot": insert a new date column (29-dec) before the existing
# "01-oct." column, shifting all the October..(end) columns one to the
# right. The newly inserted column has no data yet, so every data row gets
# the "-" placeholder (matching the other not-yet-available future dates).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

$ws.Range("EZ1").EntireColumn.Insert()

$ws.Range("EZ1").Value() = "29-dec"
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 156).Value() = "-"
}

# Sheet "Gaz": append two more rows of (not yet available) daily prices.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A184").Value() = "2025-12-27"
$wsGaz.Range("B184").Value() = ""
$wsGaz.Range("A185").Value() = "2025-12-28"
$wsGaz.Range("B185").Value() = ""
